# Auto update Excel log
#
# Appends the latest batch of sensor/alert readings (captured 2026-01-28,
# ~16:36-16:37) to the SeniorConnect master log workbook: new PIR motion
# events, Humidity/Temperature bathroom sensor samples, and the ALERTS
# entry generated when the resident entered the bathroom.

$wb = $excel.ActiveWorkbook

# --- "PIR" sheet ---
$ws = $wb.Worksheets.Item("PIR")
$pirRows = New-Object System.Collections.ArrayList
$null = $pirRows.Add(@("2026-01-28", "16:36:45", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:36:47", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:36:52", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:36:54", "16:00", "Bathroom", "Motion Detected", "Active"))
$null = $pirRows.Add(@("2026-01-28", "16:37:01", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:06", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:11", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:16", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:21", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:26", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:32", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:37", "16:00", "Bathroom", "No Motion", "Inactive"))
$null = $pirRows.Add(@("2026-01-28", "16:37:42", "16:00", "Bathroom", "No Motion", "Inactive"))
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$endRow = $startRow + $pirRows.Count - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "@"
$r = $startRow
foreach ($row in $pirRows) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# --- "Humidity" sheet ---
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = New-Object System.Collections.ArrayList
$null = $humidityRows.Add(@("2026-01-28", "16:36:45", "16:00", "Bathroom", "87.1%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:36:54", "16:00", "Bathroom", "87.1%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:36:58", "16:00", "Bathroom", "88.0%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:02", "16:00", "Bathroom", "87.0%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:06", "16:00", "Bathroom", "87.9%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:14", "16:00", "Bathroom", "87.0%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:18", "16:00", "Bathroom", "87.9%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:22", "16:00", "Bathroom", "86.9%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:26", "16:00", "Bathroom", "87.9%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:34", "16:00", "Bathroom", "87.8%", "Active"))
$null = $humidityRows.Add(@("2026-01-28", "16:37:42", "16:00", "Bathroom", "87.0%", "Active"))
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$endRow = $startRow + $humidityRows.Count - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "@"
$ws.Range($ws.Cells.Item($startRow, 5), $ws.Cells.Item($endRow, 5)).NumberFormat = "@"
$r = $startRow
foreach ($row in $humidityRows) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# --- "Temperature" sheet ---
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = New-Object System.Collections.ArrayList
$null = $temperatureRows.Add(@("2026-01-28", "16:36:45", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:36:54", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:36:58", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:02", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:06", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:14", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:18", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:22", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:26", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:34", "16:00", "Bathroom", "22.8C", "Active"))
$null = $temperatureRows.Add(@("2026-01-28", "16:37:42", "16:00", "Bathroom", "22.9C", "Active"))
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$endRow = $startRow + $temperatureRows.Count - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "@"
$r = $startRow
foreach ($row in $temperatureRows) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# --- "ALERTS" sheet ---
$ws = $wb.Worksheets.Item("ALERTS")
$alertsRows = New-Object System.Collections.ArrayList
$null = $alertsRows.Add(@("2026-01-28", "16:37:39", "16:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom"))
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$endRow = $startRow + $alertsRows.Count - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).NumberFormat = "@"
$r = $startRow
foreach ($row in $alertsRows) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}
